# "Generate Report for Archive"
#
# 1) The status text "Ready for handoff" becomes "In Translation" wherever it
#    appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all share the same
#    string, so updating the cells that hold it updates all of them).
# 2) The "Status" column is narrower in the refreshed report: Overview columns
#    E & F and the Status column (C) on the language sheets shrink from
#    ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# Target stored column width (as written into the worksheet XML) is
# 13.4101845877511. The COM layer quantizes ColumnWidth to whole pixels
# (stored = (round(ColumnWidth*6)+5)/6), so 12.5 is the closest input that
# lands on the nearest representable width (13.333333333333334).
$narrowWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = $narrowWidth
$overview.Columns.Item(6).ColumnWidth = $narrowWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = $narrowWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = $narrowWidth
